# Applies the two textual changes described by the diff:
#  1. Slide 1 "date" placeholder: November 17, 2025 -> November 18, 2025
#  2. Slide 5 phase/timeline placeholder: remove the Markdown-style
#     asterisk emphasis around the "(Months N-M)" qualifiers for all
#     three phases, e.g. " *(Months 1-2)*" -> " (Months 1-2)".

function Replace-FirstOccurrence {
    param($TextRange, $OldValue, $NewValue)

    $full = $TextRange.Text
    $idx = $full.IndexOf($OldValue)
    if ($idx -ge 0) {
        $sub = $TextRange.Characters($idx + 1, $OldValue.Length)
        $sub.Text = $NewValue
    }
}

$p = $ppt.ActivePresentation

# --- Slide 1: update the date in the "Author | Date" subtitle ---
# (This whole paragraph is a single run. Re-assigning TextRange.Text
# directly injects a default <a:rPr lang="en-US"/> into the run, so
# instead we select the run's exact full span via Characters() and set
# the text on that sub-range, which keeps the original run untouched
# other than its text content.)
$slide1 = $p.Slides.Item(1)
$dateShape = $slide1.Shapes.Item(5)
$dateRange = $dateShape.TextFrame.TextRange
$dateFull = $dateRange.Characters(1, $dateRange.Text.Length)
$dateFull.Text = "Alison Smith | November 18, 2025"

# --- Slide 5: strip the asterisk emphasis from the phase month ranges ---
$slide5 = $p.Slides.Item(5)
$phaseShape = $slide5.Shapes.Item(3)
$phaseRange = $phaseShape.TextFrame.TextRange

Replace-FirstOccurrence $phaseRange " *(Months 1-2)*" " (Months 1-2)"
Replace-FirstOccurrence $phaseRange " *(Months 3-4)*" " (Months 3-4)"
Replace-FirstOccurrence $phaseRange " *(Months 5-6)*" " (Months 5-6)"
